# codeforIATI SectorGroup.xlsx update:
# The "category-name" column was moved (in the shared-string table) to sit
# immediately after "category-code", which has the visible effect of
# rotating the three columns E (was group-code), F (was group-name) and
# G (was category-name) one step to the left:
#   new E (category-name) = old G
#   new F (group-code)    = old E
#   new G (group-name)    = old F
# This holds for the header row too (row 1), and for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

# First, read every row's current E/F/G values into memory so writes to
# one row can never be read back for another (not actually an issue here
# since each row is independent, but it keeps the read/write phases
# cleanly separated and easy to reason about).
$oldE = @()
$oldF = @()
$oldG = @()

for ($r = 1; $r -le $rowCount; $r++) {
    $oldE += , ($ws.Cells.Item($r, 5).Value())
    $oldF += , ($ws.Cells.Item($r, 6).Value())
    $oldG += , ($ws.Cells.Item($r, 7).Value())
}

for ($r = 1; $r -le $rowCount; $r++) {
    $i = $r - 1
    $ws.Cells.Item($r, 5).Value = $oldG[$i]
    $ws.Cells.Item($r, 6).Value = $oldE[$i]
    $ws.Cells.Item($r, 7).Value = $oldF[$i]
}
